$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the average mark and subject values for the remaining rows
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = "Assembler"

$ws.Range("C3").Value = "Programming"

$ws.Range("C4").Value = "WinApi"

# The last two test rows (Winter/2020 Programming, Winter/2019 YAPVU) are no
# longer part of the table, so remove them entirely.
$ws.Range("A5:C6").Delete()
